# "add hospital table script." — append the `hospital` table definition
# to the data dictionary sheet (rows 50-58), two blank rows below the
# existing `serverroleoperationmatrix` table (which ends at row 47),
# matching the two-blank-row convention already used between every other
# table block on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: hospitalid (PK, int, not null)
$ws.Range("A50").Value = "hospital"
$ws.Range("B50").Value = "hospitalid"
$ws.Range("C50").Value = "int"
$ws.Range("E50").Value = "no"
$ws.Range("F50").Value = "pk"

# Row 51: hospitalcode (UK, varchar(10), not null)
$ws.Range("A51").Value = "hospital"
$ws.Range("B51").Value = "hospitalcode"
$ws.Range("C51").Value = "varchar"
$ws.Range("D51").Value = 10
$ws.Range("E51").Value = "no"
$ws.Range("F51").Value = "uk"

# Row 52: hospital name (varchar(255), not null)
$ws.Range("A52").Value = "hospital "
$ws.Range("B52").Value = "hospital name"
$ws.Range("C52").Value = "varchar"
$ws.Range("D52").Value = 255
$ws.Range("E52").Value = "no"

# Row 53: shortname (UK, varchar(10), nullable)
$ws.Range("A53").Value = "hospital "
$ws.Range("B53").Value = "shortname"
$ws.Range("C53").Value = "varchar"
$ws.Range("D53").Value = 10
$ws.Range("E53").Value = "yes"
$ws.Range("F53").Value = "uk"

# Row 54: address (varchar, nullable)
$ws.Range("A54").Value = "hospital"
$ws.Range("B54").Value = "address"
$ws.Range("C54").Value = "varchar"
$ws.Range("E54").Value = "yes"

# Row 55: phone (varchar(15), nullable)
$ws.Range("A55").Value = "hospital"
$ws.Range("B55").Value = "phone"
$ws.Range("C55").Value = "varchar"
$ws.Range("D55").Value = 15
$ws.Range("E55").Value = "yes"

# Row 56: databasename (varchar(10), not null)
$ws.Range("A56").Value = "hospital"
$ws.Range("B56").Value = "databasename"
$ws.Range("C56").Value = "varchar"
$ws.Range("D56").Value = 10
$ws.Range("E56").Value = "no"

# Row 57: connectionstring (varchar(max), not null)
$ws.Range("A57").Value = "hospital"
$ws.Range("B57").Value = "connectionstring"
$ws.Range("C57").Value = "varchar(max)"
$ws.Range("D57").Value = "custome"
$ws.Range("E57").Value = "no"

# Row 58: hospitallogo (varchar, nullable)
$ws.Range("A58").Value = "hospital"
$ws.Range("B58").Value = "hospitallogo"
$ws.Range("C58").Value = "varchar"
$ws.Range("E58").Value = "yes"

# Leave the selection where the author's cursor ended up after typing the
# new table in (matches the committed sheetView selection).
[void]$ws.Range("G52").Select()
